$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.618.67'
$ws.Range("E2").Value = '  +0.91%  '
$ws.Range("D3").Value = '1.630.69'
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''212.70'
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D6").Value = '''0.498'
$ws.Range("E6").Value = '  +2.65%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  +1.48%  '
$ws.Range("E9").Value = '  +0.86%  '
$ws.Range("D10").Value = '''19.13'
$ws.Range("E10").Value = '  +1.44%  '
$ws.Range("E11").Value = '  +3.36%  '
$ws.Range("D12").Value = '1.858.14'
$ws.Range("D13").Value = '1.603.87'
$ws.Range("E13").Value = '  -1.24%  '
$ws.Range("D14").Value = '''4.09'
$ws.Range("E14").Value = '  +1.77%  '
$ws.Range("E15").Value = '  +0.92%  '
$ws.Range("D16").Value = '26.612.66'
$ws.Range("E16").Value = '  +0.95%  '
$ws.Range("D17").Value = '''63.21'
$ws.Range("E17").Value = '  +1.20%  '
$ws.Range("D18").Value = '0.0₃0740'
$ws.Range("E18").Value = '  +1.76%  '
$ws.Range("D19").Value = '''217.85'
$ws.Range("E19").Value = '  +7.61%  '
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("E21").Value = '  +0.27%  '
$ws.Range("E22").Value = '  +2.01%  '
$ws.Range("D23").Value = '''9.33'
$ws.Range("E24").Value = '  +3.51%  '
$ws.Range("E25").Value = '  +2.39%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("E27").Value = '  +0.73%  '
$ws.Range("D28").Value = '''6.83'
$ws.Range("E28").Value = '  +3.83%  '
$ws.Range("D29").Value = '''15.52'
$ws.Range("E29").Value = '  +2.14%  '
$ws.Range("D30").Value = '''0.0503'
$ws.Range("E30").Value = '  -3.45%  '
$ws.Range("D32").Value = '''3.30'
$ws.Range("E32").Value = '  +3.70%  '
$ws.Range("E33").Value = '  +1.92%  '
$ws.Range("E34").Value = '  +0.33%  '
$ws.Range("D35").Value = '''2.40'
$ws.Range("E35").Value = '  +0.29%  '
$ws.Range("D36").Value = '1.213.65'
$ws.Range("E36").Value = '  +3.16%  '
$ws.Range("D37").Value = '''0.0171'
$ws.Range("E37").Value = '  +4.58%  '
$ws.Range("D38").Value = '''0.806'
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("E39").Value = '  +0.12%  '
$ws.Range("D40").Value = '''0.501'
$ws.Range("E40").Value = '  +0.70%  '
$ws.Range("E41").Value = '  -1.92%  '
$ws.Range("E42").Value = '  +0.35%  '
$ws.Range("D43").Value = '''0.790'
$ws.Range("E43").Value = '  +0.57%  '
$ws.Range("D44").Value = '1.772.00'
$ws.Range("E44").Value = '  +0.67%  '
$ws.Range("D45").Value = '''92.69'
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("E46").Value = '  +2.00%  '
$ws.Range("E47").Value = '  +1.62%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0103'
$ws.Range("E48").Value = '  -2.09%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.0511'
$ws.Range("E49").Value = '  +0.52%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''7.60'
$ws.Range("E50").Value = '  +4.15%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '''0.410'
$ws.Range("E51").Value = '  -0.08%  '
